# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" sheet right after "总计", shifting "2022-Q2",
# "2022-Q1" and "2021-Q4" one position to the right, and adds the
# corresponding summary row to "总计".
#
# The host's worksheet handles are resolved *by position* every time a
# property/method is touched (not by stable object identity), so any
# previously-grabbed handle to "sheet index N" silently starts pointing at
# whatever sheet now lives at index N after an Add/Delete. To stay correct
# we rebuild the whole tab strip from scratch in final left-to-right order,
# finishing (name + content) each sheet immediately after creating it and
# never reusing a handle across a structural change.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Set-PlainText($ws, $cell, $val) {
    # Force text typing (even for numeric-looking strings like "161123" or
    # "4.78") without leaving a stray NumberFormat-driven style behind.
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($cell).ClearFormats()
}

function Set-StyledText($ws, $cell, $val, $styleSrc) {
    Set-PlainText $ws $cell $val
    $styleSrc.Copy()
    $ws.Range($cell).PasteSpecial($xlPasteFormats)
}

function Set-StyledNumber($ws, $cell, $val, $styleSrc) {
    $styleSrc.Copy()
    $ws.Range($cell).PasteSpecial($xlPasteFormats)
    $ws.Range($cell).Value = $val
}

# ---------------------------------------------------------------------
# 1. Tear down every sheet except "总计" (kept as the style/style-id
#    anchor and as sheet index 1 throughout).
# ---------------------------------------------------------------------
$wb.Worksheets.Item(4).Delete() # 2021-Q4
$wb.Worksheets.Item(3).Delete() # 2022-Q1
$wb.Worksheets.Item(2).Delete() # 2022-Q2

$total = $wb.Worksheets.Item(1)
$headerStyleSrc = $total.Range("B1")
$colAStyleSrc = $total.Range("A2")

# ---------------------------------------------------------------------
# 2. Rebuild "总计" rows 2-5 (row 1 header is untouched).
# ---------------------------------------------------------------------
Set-StyledNumber $total "A2" 0 $colAStyleSrc
Set-PlainText $total "B2" "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.13

Set-StyledNumber $total "A3" 1 $colAStyleSrc
Set-PlainText $total "B3" "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.21

Set-StyledNumber $total "A4" 2 $colAStyleSrc
Set-PlainText $total "B4" "2022-Q1"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.17

Set-StyledNumber $total "A5" 3 $colAStyleSrc
Set-PlainText $total "B5" "2021-Q4"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.3

# ---------------------------------------------------------------------
# 3. Re-add the four quarter sheets in final order, populating each one
#    completely before the next insertion shifts it out from under us.
# ---------------------------------------------------------------------

# --- 2022-Q3 (new) ---
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

Set-StyledText $q3 "B1" "基金代码" $headerStyleSrc
Set-StyledText $q3 "C1" "基金名称" $headerStyleSrc
Set-StyledText $q3 "D1" "基金规模" $headerStyleSrc
Set-StyledText $q3 "E1" "股票总仓位" $headerStyleSrc
Set-StyledText $q3 "F1" "仓位占比" $headerStyleSrc
Set-StyledText $q3 "G1" "持有市值(亿元)" $headerStyleSrc
Set-StyledText $q3 "H1" "仓位排名" $headerStyleSrc

Set-StyledNumber $q3 "A2" 0 $colAStyleSrc
Set-PlainText $q3 "B2" "161123"
Set-PlainText $q3 "C2" "易方达并购重组指数（LOF）"
Set-PlainText $q3 "D2" "4.34"
Set-PlainText $q3 "E2" "94.11"
Set-PlainText $q3 "F2" "3.04"
Set-PlainText $q3 "G2" "0.1319"
$q3.Range("H2").Value = 6

# --- 2022-Q2 (unchanged data, shifted one slot right) ---
$q2 = $wb.Worksheets.Add($null, $q3)
$q2.Name = "2022-Q2"

Set-StyledText $q2 "B1" "基金代码" $headerStyleSrc
Set-StyledText $q2 "C1" "基金名称" $headerStyleSrc
Set-StyledText $q2 "D1" "基金规模" $headerStyleSrc
Set-StyledText $q2 "E1" "股票总仓位" $headerStyleSrc
Set-StyledText $q2 "F1" "仓位占比" $headerStyleSrc
Set-StyledText $q2 "G1" "持有市值(亿元)" $headerStyleSrc
Set-StyledText $q2 "H1" "仓位排名" $headerStyleSrc

Set-StyledNumber $q2 "A2" 0 $colAStyleSrc
Set-PlainText $q2 "B2" "161123"
Set-PlainText $q2 "C2" "易方达并购重组指数（LOF）"
Set-PlainText $q2 "D2" "4.78"
Set-PlainText $q2 "E2" "94.61"
Set-PlainText $q2 "F2" "2.70"
Set-PlainText $q2 "G2" "0.1291"
$q2.Range("H2").Value = 9

Set-StyledNumber $q2 "A3" 1 $colAStyleSrc
Set-PlainText $q2 "B3" "001050"
Set-PlainText $q2 "C3" "汇添富成长多因子量化策略股票"
Set-PlainText $q2 "D3" "8.24"
Set-PlainText $q2 "E3" "93.77"
Set-PlainText $q2 "F3" "0.96"
Set-PlainText $q2 "G3" "0.0791"
$q2.Range("H3").Value = 10

# --- 2022-Q1 (unchanged data, shifted one slot right) ---
$q1 = $wb.Worksheets.Add($null, $q2)
$q1.Name = "2022-Q1"

Set-StyledText $q1 "B1" "基金代码" $headerStyleSrc
Set-StyledText $q1 "C1" "基金名称" $headerStyleSrc
Set-StyledText $q1 "D1" "基金规模" $headerStyleSrc
Set-StyledText $q1 "E1" "股票总仓位" $headerStyleSrc
Set-StyledText $q1 "F1" "仓位占比" $headerStyleSrc
Set-StyledText $q1 "G1" "持有市值(亿元)" $headerStyleSrc
Set-StyledText $q1 "H1" "仓位排名" $headerStyleSrc

Set-StyledNumber $q1 "A2" 0 $colAStyleSrc
Set-PlainText $q1 "B2" "161123"
Set-PlainText $q1 "C2" "易方达并购重组指数（LOF）"
Set-PlainText $q1 "D2" "4.78"
Set-PlainText $q1 "E2" "94.71"
Set-PlainText $q1 "F2" "3.62"
Set-PlainText $q1 "G2" "0.1730"
$q1.Range("H2").Value = 8

# --- 2021-Q4 (unchanged data, shifted one slot right) ---
$q4 = $wb.Worksheets.Add($null, $q1)
$q4.Name = "2021-Q4"

Set-StyledText $q4 "B1" "基金代码" $headerStyleSrc
Set-StyledText $q4 "C1" "基金名称" $headerStyleSrc
Set-StyledText $q4 "D1" "基金规模" $headerStyleSrc
Set-StyledText $q4 "E1" "股票总仓位" $headerStyleSrc
Set-StyledText $q4 "F1" "仓位占比" $headerStyleSrc
Set-StyledText $q4 "G1" "持有市值(亿元)" $headerStyleSrc
Set-StyledText $q4 "H1" "仓位排名" $headerStyleSrc

Set-StyledNumber $q4 "A2" 0 $colAStyleSrc
Set-PlainText $q4 "B2" "460002"
Set-PlainText $q4 "C2" "华泰柏瑞积极成长混合A"
Set-PlainText $q4 "D2" "7.95"
Set-PlainText $q4 "E2" "71.27"
Set-PlainText $q4 "F2" "2.72"
Set-PlainText $q4 "G2" "0.2162"
$q4.Range("H2").Value = 10

Set-StyledNumber $q4 "A3" 1 $colAStyleSrc
Set-PlainText $q4 "B3" "004223"
Set-PlainText $q4 "C3" "金信多策略精选灵活配置混合"
Set-PlainText $q4 "D3" "2.00"
Set-PlainText $q4 "E3" "50.19"
Set-PlainText $q4 "F3" "3.94"
Set-PlainText $q4 "G3" "0.0788"
$q4.Range("H3").Value = 4

Set-StyledNumber $q4 "A4" 2 $colAStyleSrc
Set-PlainText $q4 "B4" "960030"
Set-PlainText $q4 "C4" "华泰柏瑞积极成长混合H"
Set-PlainText $q4 "D4" "0.00"
Set-PlainText $q4 "E4" "71.27"
Set-PlainText $q4 "F4" "2.72"
$q4.Range("G4").Value = 0
$q4.Range("H4").Value = 10
